# SCD0027 - Mengunci Gembok
# Regression-run update: bump the TEXT1 run counter (column L) on rows 2 and 4
# from 1121 to 1177, and record the generated FILE2 (column Q) zip name on
# row 3, then leave the selection on the last-touched cell (L2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TEXT1 column (L) run id bump: 1121 -> 1177
$ws.Range("L2").Value = 1177
$ws.Range("L4").Value = 1177

# FILE2 column (Q) on row 3: newly recorded downloaded/uploaded file name
$ws.Range("Q3").Value = "Wilayah_05_202206_31.zip"

# Match the author's final selection (cell L2) left active in the sheet view
$ws.Range("L2").Select()
